# Scheduled market-data refresh: update currentAveragePrice / LevePrice /
# LeveProfit columns (H:N) across the per-job Leve tables (ALC, ARM, BSM,
# CRP, CUL, GSM, LTW, WVR) with freshly pulled Universalis price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 932.8333
$ws.Range("I19").Value = 932.8333
$ws.Range("K19").Value = 932.8333
$ws.Range("M19").Value = -757.8333

# Row 33
$ws.Range("H33").Value = 2502955.2
$ws.Range("I33").Value = 2859814.2
$ws.Range("J33").Value = 4942
$ws.Range("K33").Value = 2859814.2
$ws.Range("L33").Value = 4942
$ws.Range("M33").Value = -2859585.2
$ws.Range("N33").Value = -5400

# Row 132
$ws.Range("H132").Value = 3576.7144
$ws.Range("I132").Value = 3387.0557
$ws.Range("J132").Value = 4714.6665
$ws.Range("K132").Value = 10161.1671
$ws.Range("L132").Value = 14143.9995
$ws.Range("M132").Value = -7631.167099999999
$ws.Range("N132").Value = -19203.9995

# Row 137
$ws.Range("H137").Value = 839344.4399999999
$ws.Range("I137").Value = 629994.75
$ws.Range("J137").Value = 1258043.8
$ws.Range("K137").Value = 1889984.25
$ws.Range("L137").Value = 3774131.4
$ws.Range("M137").Value = -1887434.25
$ws.Range("N137").Value = -3779231.4

# Row 138
$ws.Range("H138").Value = 5386.7627
$ws.Range("I138").Value = 2890.8076
$ws.Range("J138").Value = 6588.5186
$ws.Range("K138").Value = 8672.4228
$ws.Range("L138").Value = 19765.5558
$ws.Range("M138").Value = -3532.4228
$ws.Range("N138").Value = -30045.5558

# Row 141
$ws.Range("H141").Value = 4293.769
$ws.Range("I141").Value = 2646.5557
$ws.Range("J141").Value = 8000
$ws.Range("K141").Value = 7939.6671
$ws.Range("L141").Value = 24000
$ws.Range("M141").Value = -2759.6671
$ws.Range("N141").Value = -34360

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2716.5432
$ws.Range("I32").Value = 1785.3164
$ws.Range("K32").Value = 1785.3164
$ws.Range("M32").Value = -1498.3164

# Row 61
$ws.Range("H61").Value = 4133.222
$ws.Range("I61").Value = 3837.375
$ws.Range("K61").Value = 3837.375
$ws.Range("M61").Value = -3625.375

# Row 107
$ws.Range("H107").Value = 50000
$ws.Range("J107").Value = 50000
$ws.Range("L107").Value = 50000
$ws.Range("N107").Value = -57680

# Row 132
$ws.Range("H132").Value = 8399.799999999999
$ws.Range("I132").Value = 4999.5
$ws.Range("J132").Value = 10666.667
$ws.Range("K132").Value = 14998.5
$ws.Range("L132").Value = 32000.001
$ws.Range("M132").Value = -12468.5
$ws.Range("N132").Value = -37060.001

# Row 135
$ws.Range("H135").Value = 74285
$ws.Range("J135").Value = 74285
$ws.Range("L135").Value = 74285
$ws.Range("N135").Value = -84425

# Row 136
$ws.Range("H136").Value = 4133.222
$ws.Range("I136").Value = 3837.375
$ws.Range("K136").Value = 11512.125
$ws.Range("M136").Value = -8962.125

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 200003920
$ws.Range("I20").Value = 500002300
$ws.Range("J20").Value = 5004.6665
$ws.Range("K20").Value = 500002300
$ws.Range("L20").Value = 5004.6665
$ws.Range("M20").Value = -500002053
$ws.Range("N20").Value = -5498.6665

# Row 134
$ws.Range("H134").Value = 46166
$ws.Range("I134").Value = 4453.8184
$ws.Range("J134").Value = 505000
$ws.Range("K134").Value = 13361.4552
$ws.Range("L134").Value = 1515000
$ws.Range("M134").Value = -10826.4552
$ws.Range("N134").Value = -1520070

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 75647
$ws.Range("I31").Value = 2116
$ws.Range("J31").Value = 208002.8
$ws.Range("K31").Value = 2116
$ws.Range("L31").Value = 208002.8
$ws.Range("M31").Value = -1821
$ws.Range("N31").Value = -208592.8

# Row 34
$ws.Range("H34").Value = 75647
$ws.Range("I34").Value = 2116
$ws.Range("J34").Value = 208002.8
$ws.Range("K34").Value = 2116
$ws.Range("L34").Value = 208002.8
$ws.Range("M34").Value = -1914
$ws.Range("N34").Value = -208406.8

# Row 58
$ws.Range("H58").Value = 166811.53
$ws.Range("I58").Value = 198232.44
$ws.Range("J58").Value = 6564.9
$ws.Range("K58").Value = 198232.44
$ws.Range("L58").Value = 6564.9
$ws.Range("M58").Value = -198029.44
$ws.Range("N58").Value = -6970.9

# Row 86
$ws.Range("H86").Value = 41763.2
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 41763.2
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 41763.2
$ws.Range("N86").Value = -44009.2
$ws.Range("M86").ClearContents()

# Row 89
$ws.Range("H89").Value = 41763.2
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 41763.2
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 208816
$ws.Range("N89").Value = -220048
$ws.Range("M89").ClearContents()

# Row 99
$ws.Range("H99").Value = 364413.94
$ws.Range("I99").Value = 8466.333000000001
$ws.Range("J99").Value = 631374.6
$ws.Range("K99").Value = 8466.333000000001
$ws.Range("L99").Value = 631374.6
$ws.Range("M99").Value = -6968.333000000001
$ws.Range("N99").Value = -634370.6

# Row 107
$ws.Range("H107").Value = 306.83334
$ws.Range("I107").Value = 193.2
$ws.Range("J107").Value = 875
$ws.Range("K107").Value = 193.2
$ws.Range("L107").Value = 875
$ws.Range("M107").Value = 1726.8
$ws.Range("N107").Value = -4715

# Row 126
$ws.Range("H126").Value = 364413.94
$ws.Range("I126").Value = 8466.333000000001
$ws.Range("J126").Value = 631374.6
$ws.Range("K126").Value = 25398.999
$ws.Range("L126").Value = 1894123.8
$ws.Range("M126").Value = -22928.999
$ws.Range("N126").Value = -1899063.8

# Row 132
$ws.Range("H132").Value = 3004.487
$ws.Range("I132").Value = 2684.2903
$ws.Range("J132").Value = 4245.25
$ws.Range("K132").Value = 8052.8709
$ws.Range("L132").Value = 12735.75
$ws.Range("M132").Value = -5522.8709
$ws.Range("N132").Value = -17795.75

# Row 134
$ws.Range("H134").Value = 1543807.8
$ws.Range("I134").Value = 839125
$ws.Range("J134").Value = 10000000
$ws.Range("K134").Value = 2517375
$ws.Range("L134").Value = 30000000
$ws.Range("M134").Value = -2514840
$ws.Range("N134").Value = -30005070

# Row 136
$ws.Range("H136").Value = 166811.53
$ws.Range("I136").Value = 198232.44
$ws.Range("J136").Value = 6564.9
$ws.Range("K136").Value = 594697.3200000001
$ws.Range("L136").Value = 19694.7
$ws.Range("M136").Value = -592147.3200000001
$ws.Range("N136").Value = -24794.7

$ws = $wb.Worksheets.Item("CUL")
# Row 132
$ws.Range("H132").Value = 481987.1
$ws.Range("J132").Value = 838083.5
$ws.Range("L132").Value = 7542751.5
$ws.Range("N132").Value = -7547811.5

# Row 139
$ws.Range("H139").Value = 3206.276
$ws.Range("I139").Value = 1681.4546
$ws.Range("J139").Value = 7998.5713
$ws.Range("K139").Value = 5044.3638
$ws.Range("L139").Value = 23995.7139
$ws.Range("M139").Value = 95.63619999999992
$ws.Range("N139").Value = -34275.7139

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 8481.799999999999
$ws.Range("I70").Value = 10669.667
$ws.Range("K70").Value = 10669.667
$ws.Range("M70").Value = -10399.667

# Row 73
$ws.Range("H73").Value = 8481.799999999999
$ws.Range("I73").Value = 10669.667
$ws.Range("K73").Value = 10669.667
$ws.Range("M73").Value = -9733.666999999999

# Row 80
$ws.Range("H80").Value = 1671267.4
$ws.Range("I80").Value = 2005260.6
$ws.Range("J80").Value = 1432700.9
$ws.Range("K80").Value = 2005260.6
$ws.Range("L80").Value = 1432700.9
$ws.Range("M80").Value = -2004262.6
$ws.Range("N80").Value = -1434696.9

# Row 83
$ws.Range("H83").Value = 1671267.4
$ws.Range("I83").Value = 2005260.6
$ws.Range("J83").Value = 1432700.9
$ws.Range("K83").Value = 10026303
$ws.Range("L83").Value = 7163504.5
$ws.Range("M83").Value = -10021311
$ws.Range("N83").Value = -7173488.5

# Row 126
$ws.Range("H126").Value = 2746.3845
$ws.Range("I126").Value = 2070.3
$ws.Range("K126").Value = 6210.900000000001
$ws.Range("M126").Value = -3740.900000000001

# Row 132
$ws.Range("H132").Value = 270770.62
$ws.Range("I132").Value = 288330.5
$ws.Range("K132").Value = 864991.5
$ws.Range("M132").Value = -862461.5

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()

# Row 132
$ws.Range("H132").Value = 2621.2122
$ws.Range("I132").Value = 1546.08
$ws.Range("K132").Value = 4638.24
$ws.Range("M132").Value = -2108.24

# Row 136
$ws.Range("H136").Value = 402769.9
$ws.Range("I136").Value = 456783.97
$ws.Range("J136").Value = 6667
$ws.Range("K136").Value = 1370351.91
$ws.Range("L136").Value = 20001
$ws.Range("M136").Value = -1367801.91
$ws.Range("N136").Value = -25101

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 49430.91
$ws.Range("I132").Value = 2304.5
$ws.Range("J132").Value = 261499.75
$ws.Range("K132").Value = 6913.5
$ws.Range("L132").Value = 784499.25
$ws.Range("M132").Value = -4383.5
$ws.Range("N132").Value = -789559.25

